$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.255781
$ws.Range("H2").Value = 3.767343
$ws.Range("I2").Value = 0.01633546530699055
$ws.Range("J2").Value = 0.01633546530699055
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.460162333333333
$ws.Range("N2").Value = 28.380487
$ws.Range("O2").Value = 0.08197024919772995
$ws.Range("P2").Value = 0.08197024919772995
$ws.Range("Q2").Value = 11.87989211511567
$ws.Range("R2").Value = 106.919029036041
$ws.Range("S2").Value = 0.001339022161974888
$ws.Range("T2").Value = 0.001339022161974888
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.255781
$ws.Range("H3").Value = 3.767343
$ws.Range("I3").Value = 0.01633546530699055
$ws.Range("J3").Value = 0.01633546530699055
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.18256633333333
$ws.Range("N3").Value = 138.547699
$ws.Range("O3").Value = 0.4001618933742075
$ws.Range("P3").Value = 0.4001618933742075
$ws.Range("Q3").Value = 57.99518933263967
$ws.Range("R3").Value = 521.9567039937571
$ws.Range("S3").Value = 0.006536830726394018
$ws.Range("T3").Value = 0.006536830726394018
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.255781
$ws.Range("H4").Value = 3.767343
$ws.Range("I4").Value = 0.01633546530699055
$ws.Range("J4").Value = 0.01633546530699055
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 51.40166966666666
$ws.Range("N4").Value = 154.205009
$ws.Range("O4").Value = 0.4453842886934318
$ws.Range("P4").Value = 0.4453842886934319
$ws.Range("Q4").Value = 64.54924013567633
$ws.Range("R4").Value = 580.9431612210871
$ws.Range("S4").Value = 0.00727555959623022
$ws.Range("T4").Value = 0.007275559596230221
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.255781
$ws.Range("H5").Value = 3.767343
$ws.Range("I5").Value = 0.01633546530699055
$ws.Range("J5").Value = 0.01633546530699055
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.365307333333334
$ws.Range("N5").Value = 25.095922
$ws.Range("O5").Value = 0.07248356873463072
$ws.Range("P5").Value = 0.07248356873463073
$ws.Range("Q5").Value = 10.50499400836067
$ws.Range("R5").Value = 94.54494607524602
$ws.Range("S5").Value = 0.001184052822391425
$ws.Range("T5").Value = 0.001184052822391425
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 64.17341733333332
$ws.Range("H6").Value = 192.520252
$ws.Range("I6").Value = 0.8347814089237634
$ws.Range("J6").Value = 0.8347814089237634
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.460162333333333
$ws.Range("N6").Value = 28.380487
$ws.Range("O6").Value = 0.08197024919772995
$ws.Range("P6").Value = 0.08197024919772995
$ws.Range("Q6").Value = 607.0909454580802
$ws.Range("R6").Value = 5463.818509122723
$ws.Range("S6").Value = 0.068427240115113
$ws.Range("T6").Value = 0.068427240115113
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 64.17341733333332
$ws.Range("H7").Value = 192.520252
$ws.Range("I7").Value = 0.8347814089237634
$ws.Range("J7").Value = 0.8347814089237634
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.18256633333333
$ws.Range("N7").Value = 138.547699
$ws.Range("O7").Value = 0.4001618933742075
$ws.Range("P7").Value = 0.4001618933742075
$ws.Range("Q7").Value = 2963.693102833349
$ws.Range("R7").Value = 26673.23792550014
$ws.Range("S7").Value = 0.3340477091485217
$ws.Range("T7").Value = 0.3340477091485217
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 64.17341733333332
$ws.Range("H8").Value = 192.520252
$ws.Range("I8").Value = 0.8347814089237634
$ws.Range("J8").Value = 0.8347814089237634
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 51.40166966666666
$ws.Range("N8").Value = 154.205009
$ws.Range("O8").Value = 0.4453842886934318
$ws.Range("P8").Value = 0.4453842886934319
$ws.Range("Q8").Value = 3298.62079914914
$ws.Range("R8").Value = 29687.58719234226
$ws.Range("S8").Value = 0.3717985240280112
$ws.Range("T8").Value = 0.3717985240280113
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 64.17341733333332
$ws.Range("H9").Value = 192.520252
$ws.Range("I9").Value = 0.8347814089237634
$ws.Range("J9").Value = 0.8347814089237634
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.365307333333334
$ws.Range("N9").Value = 25.095922
$ws.Range("O9").Value = 0.07248356873463072
$ws.Range("P9").Value = 0.07248356873463073
$ws.Range("Q9").Value = 536.8303586235937
$ws.Range("R9").Value = 4831.473227612343
$ws.Range("S9").Value = 0.06050793563211748
$ws.Range("T9").Value = 0.06050793563211749
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.135113
$ws.Range("H10").Value = 6.405339000000001
$ws.Range("I10").Value = 0.02777400226472969
$ws.Range("J10").Value = 0.02777400226472969
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.460162333333333
$ws.Range("N10").Value = 28.380487
$ws.Range("O10").Value = 0.08197024919772995
$ws.Range("P10").Value = 0.08197024919772995
$ws.Range("Q10").Value = 20.19851558001033
$ws.Range("R10").Value = 181.786640220093
$ws.Range("S10").Value = 0.002276641886858209
$ws.Range("T10").Value = 0.002276641886858209
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.135113
$ws.Range("H11").Value = 6.405339000000001
$ws.Range("I11").Value = 0.02777400226472969
$ws.Range("J11").Value = 0.02777400226472969
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.18256633333333
$ws.Range("N11").Value = 138.547699
$ws.Range("O11").Value = 0.4001618933742075
$ws.Range("P11").Value = 0.4001618933742075
$ws.Range("Q11").Value = 98.60499775166234
$ws.Range("R11").Value = 887.4449797649611
$ws.Range("S11").Value = 0.01111409733283376
$ws.Range("T11").Value = 0.01111409733283376
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.135113
$ws.Range("H12").Value = 6.405339000000001
$ws.Range("I12").Value = 0.02777400226472969
$ws.Range("J12").Value = 0.02777400226472969
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 51.40166966666666
$ws.Range("N12").Value = 154.205009
$ws.Range("O12").Value = 0.4453842886934318
$ws.Range("P12").Value = 0.4453842886934319
$ws.Range("Q12").Value = 109.7483731270057
$ws.Range("R12").Value = 987.735358143051
$ws.Range("S12").Value = 0.0123701042428464
$ws.Range("T12").Value = 0.0123701042428464
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.135113
$ws.Range("H13").Value = 6.405339000000001
$ws.Range("I13").Value = 0.02777400226472969
$ws.Range("J13").Value = 0.02777400226472969
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.365307333333334
$ws.Range("N13").Value = 25.095922
$ws.Range("O13").Value = 0.07248356873463072
$ws.Range("P13").Value = 0.07248356873463073
$ws.Range("Q13").Value = 17.86087643639533
$ws.Range("R13").Value = 160.747887927558
$ws.Range("S13").Value = 0.002013158802191324
$ws.Range("T13").Value = 0.002013158802191324
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.310205333333334
$ws.Range("H14").Value = 27.930616
$ws.Range("I14").Value = 0.1211091235045164
$ws.Range("J14").Value = 0.1211091235045164
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 9.460162333333333
$ws.Range("N14").Value = 28.380487
$ws.Range("O14").Value = 0.08197024919772995
$ws.Range("P14").Value = 0.08197024919772995
$ws.Range("Q14").Value = 88.07605380999911
$ws.Range("R14").Value = 792.6844842899919
$ws.Range("S14").Value = 0.009927345033783862
$ws.Range("T14").Value = 0.00992734503378386
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.310205333333334
$ws.Range("H15").Value = 27.930616
$ws.Range("I15").Value = 0.1211091235045164
$ws.Range("J15").Value = 0.1211091235045164
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 46.18256633333333
$ws.Range("N15").Value = 138.547699
$ws.Range("O15").Value = 0.4001618933742075
$ws.Range("P15").Value = 0.4001618933742075
$ws.Range("Q15").Value = 429.9691753836205
$ws.Range("R15").Value = 3869.722578452584
$ws.Range("S15").Value = 0.048463256166458
$ws.Range("T15").Value = 0.048463256166458
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.310205333333334
$ws.Range("H16").Value = 27.930616
$ws.Range("I16").Value = 0.1211091235045164
$ws.Range("J16").Value = 0.1211091235045164
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 51.40166966666666
$ws.Range("N16").Value = 154.205009
$ws.Range("O16").Value = 0.4453842886934318
$ws.Range("P16").Value = 0.4453842886934319
$ws.Range("Q16").Value = 478.5600990728382
$ws.Range("R16").Value = 4307.040891655543
$ws.Range("S16").Value = 0.05394010082634401
$ws.Range("T16").Value = 0.05394010082634401
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.310205333333334
$ws.Range("H17").Value = 27.930616
$ws.Range("I17").Value = 0.1211091235045164
$ws.Range("J17").Value = 0.1211091235045164
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.365307333333334
$ws.Range("N17").Value = 25.095922
$ws.Range("O17").Value = 0.07248356873463072
$ws.Range("P17").Value = 0.07248356873463073
$ws.Range("Q17").Value = 77.88272894977246
$ws.Range("R17").Value = 700.9445605479521
$ws.Range("S17").Value = 0.008778421477930494
$ws.Range("T17").Value = 0.008778421477930494
